$wb = $excel.ActiveWorkbook

# --- Create "Item sets for nobles" sheet (placed after Bretagne) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws5.Name = "Item sets for nobles"

$ws5.Range("A1").Value = "itm_h_bascinet_great,itm_heraldic_churburg_13_brass_tabard,itm_b_shynbaulds,itm_g_hourglass_gauntlets"
$ws5.Range("A2").Value = "itm_h_bascinet_oniontop,itm_heraldic_tunic_new,itm_g_finger_gauntlets,itm_b_steel_greaves_full"
$ws5.Range("A4").Value = "itm_h_klappvisier_pigface_open,itm_heraldic_churburg_13_tabard,itm_b_shynbaulds,itm_g_plate_mittens"
$ws5.Range("A5").Value = "itm_h_zitta_bascinet_novisor,itm_heraldic_early_transitional,itm_b_splinted_greaves_spurs,itm_g_wisby_gauntlets_black"
$ws5.Range("A6").Value = "itm_h_zitta_bascinet,itm_heraldic_churburg_13_brass_tabard,itm_b_shynbaulds,itm_g_hourglass_gauntlets_ornate"
$ws5.Range("A7").Value = "itm_h_klappvisier_pigface,itm_mail_long_surcoat_new_heraldic,itm_b_steel_greaves,itm_g_mail_gauntlets"
$ws5.Range("A8").Value = "itm_h_zitta_bascinet,itm_brigandine_b_heraldic,itm_b_shynbaulds,itm_g_demi_gauntlets"
$ws5.Range("A9").Value = "itm_h_zitta_bascinet_novisor,itm_heraldic_mail_tabard,itm_b_leather_boots,itm_g_mail_gauntlets"
$ws5.Range("A10").Value = "itm_h_klappvisier_pigface,itm_heraldic_early_transitional,itm_b_steel_greaves,itm_g_demi_gauntlets"
$ws5.Range("A11").Value = "itm_h_houndskull_thick,itm_heraldic_plate,itm_b_shynbaulds,itm_g_plate_mittens"
$ws5.Range("A12").Value = "itm_h_bascinet_great,itm_heraldic_tunic_new,itm_b_steel_greaves,itm_g_demi_gauntlets"
$ws5.Range("A13").Value = "itm_h_zitta_bascinet,itm_mail_long_surcoat_new_heraldic,itm_b_mail_boots,itm_g_mail_gauntlets"
$ws5.Range("A3").Value = "itm_h_hounskull_narf,itm_brigandine_b_heraldic,itm_b_shynbaulds,itm_g_hourglass_gauntlets"

$ws5.Columns.Item(1).ColumnWidth = 109.7109375
$ws5.Range("A19").Select() | Out-Null

# --- Create "Item sets for commoners" sheet (placed after "Item sets for nobles") ---
$ws6 = $wb.Worksheets.Add([Type]::Missing, $ws5)
$ws6.Name = "Item sets for commoners"

$ws6.Range("A1").Value = "Tavern Keepers"
$ws6.Range("A2").Value = "itm_a_tavern_keeper_shirt,itm_b_hosen_poulaines_custom"
$ws6.Range("A3").Value = "itm_a_commoner_apron,itm_b_hosen_shoes_custom"
$ws6.Range("A4").Value = "itm_a_peasant_coat,itm_b_ankle_boots"
$ws6.Range("A5").Value = "itm_a_peasant_shirt_white,itm_b_hosen_poulaines_custom"
$ws6.Range("A6").Value = "itm_h_felt_hat_b_black,itm_a_commoner_apron,itm_b_hosen_shoes_custom"
$ws6.Range("A7").Value = "itm_h_highlander_beret_brown,itm_a_tavern_keeper_shirt,itm_b_hosen_poulaines_custom"
$ws6.Range("A8").Value = "itm_h_highlander_beret_white,itm_a_peasant_shirt_white,itm_b_hosen_poulaines_custom"
$ws6.Range("A10").Value = "itm_a_peasant_man_custom,itm_b_hosen_poulaines_custom"
$ws6.Range("A9").Value = "itm_h_felt_hat_b_brown,itm_a_peasant_man_custom,itm_b_hosen_poulaines_custom"
$ws6.Range("A12").Value = "Misc."
$ws6.Range("A13").Value = "itm_h_felt_hat_b_white,itm_a_peasant_coat,itm_b_hosen_poulaines_custom"
$ws6.Range("A22").Value = "itm_a_tabard,itm_b_hosen_poulaines_custom"
$ws6.Range("A16").Value = "itm_h_highlander_beret_brown,itm_a_noble_shirt_brown,itm_b_hosen_poulaines_custom"
$ws6.Range("A17").Value = "itm_h_highlander_beret_red,itm_a_noble_shirt_red,itm_b_hosen_poulaines_custom"
$ws6.Range("A18").Value = "itm_h_highlander_beret_white,itm_a_noble_shirt_white,itm_b_hosen_poulaines_custom"
$ws6.Range("A19").Value = "itm_h_highlander_beret_black,itm_a_noble_shirt_black,itm_b_hosen_poulaines_custom"
$ws6.Range("A20").Value = "itm_h_highlander_beret_green,itm_a_noble_shirt_green,itm_b_hosen_poulaines_custom"
$ws6.Range("A21").Value = "itm_h_highlander_beret_blue,itm_a_noble_shirt_blue,itm_b_hosen_poulaines_custom"
$ws6.Range("A24").Value = "itm_a_commoner_apron,itm_b_leather_boots"
$ws6.Range("A25").Value = "itm_a_leather_jerkin,itm_b_leather_boots"
$ws6.Range("A26").Value = "itm_a_merchant_outfit,itm_b_hosen_shoes_custom"
$ws6.Range("A14").Value = "itm_a_peasant_man_custom,itm_b_hosen_poulaines_custom"
$ws6.Range("A15").Value = "itm_h_felt_hat_b_brown,itm_a_peasant_man_custom,itm_b_hosen_poulaines_custom"
$ws6.Range("A23").Value = "itm_a_peasant_coat,itm_b_ankle_boots"

$ws6.Columns.Item(1).ColumnWidth = 82.5703125
$ws6.Range("A26").Select() | Out-Null

# Make "Item sets for commoners" the active/selected tab, as in the target workbook
$ws6.Activate() | Out-Null
